# Refresh the cryptocurrency price/volume snapshot (cryptos.xlsx) to the
# values captured by the latest GitHub Actions run. Values that look like
# plain decimal numbers are written with a leading "'" (Excel's text-prefix)
# and then reset to the "Normal" style so they stay text cells - exactly like
# the existing (non-numeric) price strings in this sheet - instead of being
# silently re-interpreted by Excel as numbers (which would drop formatting
# such as trailing zeros, e.g. "1.00" -> 1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.053.11"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.043.92"
$ws.Range("E3").Value = "  -3.94%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'247.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.60%  "
$ws.Range("D6").Value = "'0.651"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.78%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'54.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +15.29%  "
$ws.Range("D9").Value = "'61.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").Value = "'0.373"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.57%  "
$ws.Range("D11").Value = "'0.0757"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("E12").Value = "  +4.85%  "
$ws.Range("D13").Value = "'14.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").Value = "2.338.39"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").Value = "'0.809"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.69%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'5.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.036.38"
$ws.Range("E17").Value = "  -4.27%  "
$ws.Range("D18").Value = "36.937.18"
$ws.Range("E18").Value = "  +0.01%  "
$ws.Range("D19").Value = "'71.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.19%  "
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  +5.21%  "
$ws.Range("D21").Value = "'14.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.00%  "
$ws.Range("D22").Value = "'235.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.74%  "
$ws.Range("D23").Value = "'5.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.89%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  -3.61%  "
$ws.Range("D26").Value = "'168.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.78%  "
$ws.Range("D27").Value = "'8.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("D28").Value = "'19.88"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.04%  "
$ws.Range("D29").Value = "'1.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.11%  "
$ws.Range("D30").Value = "'0.122"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("D31").Value = "'4.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").Value = "'1.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.96%  "
$ws.Range("D33").Value = "'0.0615"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("D34").Value = "'4.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.11%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'0.0871"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.71%  "
$ws.Range("D37").Value = "'2.24"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.34%  "
$ws.Range("E38").Value = "  -6.51%  "
$ws.Range("E39").Value = "  -3.28%  "
$ws.Range("D40").Value = "'0.103"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +22.29%  "
$ws.Range("D41").Value = "'18.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.71%  "
$ws.Range("B42").Value = "Gas"
$ws.Range("C42").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D42").Value = "'15.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -44.01%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.12%  "
$ws.Range("D44").Value = "'1.13"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.17%  "
$ws.Range("D45").Value = "'94.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.64%  "
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'4.06"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +41.57%  "
$ws.Range("D48").Value = "1.286.11"
$ws.Range("E48").Value = "  -5.41%  "
$ws.Range("E49").Value = "  +1.97%  "
$ws.Range("D50").Value = "'2.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("E51").Value = "  -5.98%  "
